$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "53.671.64"
$ws.Cells.Item(2,5).Value = "  -4.23%  "
$ws.Cells.Item(3,4).Value = "2.220.32"
$ws.Cells.Item(3,5).Value = "  -6.08%  "
$ws.Cells.Item(4,5).Value = "  -0.06%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "484.21"
$ws.Cells.Item(5,5).Value = "  -3.18%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "125.77"
$ws.Cells.Item(6,5).Value = "  -2.22%  "
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.999"
$ws.Cells.Item(7,5).Value = "  -0.01%  "
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.518"
$ws.Cells.Item(8,5).Value = "  -4.57%  "
$ws.Cells.Item(9,4).Value = "2.228.87"
$ws.Cells.Item(9,5).Value = "  -5.87%  "
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.0913"
$ws.Cells.Item(10,5).Value = "  -6.67%  "
$ws.Cells.Item(11,5).Value = "  -1.25%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "4.69"
$ws.Cells.Item(12,5).Value = "  -0.96%  "
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.314"
$ws.Cells.Item(13,5).Value = "  -2.50%  "
$ws.Cells.Item(14,4).Value = "2.614.34"
$ws.Cells.Item(14,5).Value = "  -6.01%  "
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "21.10"
$ws.Cells.Item(15,5).Value = "  -1.32%  "
$ws.Cells.Item(16,4).Value = "53.570.25"
$ws.Cells.Item(16,5).Value = "  -4.35%  "
$ws.Cells.Item(17,5).Value = "  -3.44%  "
$ws.Cells.Item(18,4).Value = "2.204.15"
$ws.Cells.Item(18,5).Value = "  -5.80%  "
$ws.Cells.Item(19,2).Value = "Polkadot"
$ws.Cells.Item(19,3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "3.95"
$ws.Cells.Item(19,5).Value = "  -1.92%  "
$ws.Cells.Item(20,2).Value = "Chainlink"
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "9.54"
$ws.Cells.Item(20,5).Value = "  -4.25%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "297.69"
$ws.Cells.Item(21,5).Value = "  -2.67%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "6.10"
$ws.Cells.Item(22,5).Value = "  -2.68%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "1.00"
$ws.Cells.Item(23,5).Value = "  +0.19%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "63.22"
$ws.Cells.Item(24,5).Value = "  -2.93%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "0.998"
$ws.Cells.Item(25,5).Value = "  -0.21%  "
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "0.363"
$ws.Cells.Item(26,5).Value = "  -1.16%  "
$ws.Cells.Item(27,5).Value = "  -2.42%  "
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "6.96"
$ws.Cells.Item(28,5).Value = "  -3.23%  "
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "167.88"
$ws.Cells.Item(29,5).Value = "  -2.16%  "
$ws.Cells.Item(30,5).Value = "  -3.41%  "
$ws.Cells.Item(32,4).Value = "0.0₃0673"
$ws.Cells.Item(32,5).Value = "  -5.19%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "0.997"
$ws.Cells.Item(33,5).Value = "  -0.08%  "
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "5.68"
$ws.Cells.Item(34,5).Value = "  -0.74%  "
$ws.Cells.Item(35,5).Value = "  -2.51%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "17.35"
$ws.Cells.Item(36,5).Value = "  -1.44%  "
$ws.Cells.Item(37,5).Value = "  -2.07%  "
$ws.Cells.Item(38,5).Value = "  +5.32%  "
$ws.Cells.Item(39,5).Value = "  -4.68%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "35.71"
$ws.Cells.Item(40,5).Value = "  -0.94%  "
$ws.Cells.Item(41,5).Value = "  -0.84%  "
$ws.Cells.Item(42,5).Value = "  -1.12%  "
$ws.Cells.Item(43,5).Value = "  -2.75%  "
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "122.93"
$ws.Cells.Item(44,5).Value = "  -4.91%  "
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "4.63"
$ws.Cells.Item(45,5).Value = "  -1.26%  "
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.0877"
$ws.Cells.Item(46,5).Value = "  -2.62%  "
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.533"
$ws.Cells.Item(47,5).Value = "  -5.07%  "
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "231.26"
$ws.Cells.Item(48,5).Value = "  -3.66%  "
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "0.0469"
$ws.Cells.Item(49,5).Value = "  -2.04%  "
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "0.0201"
$ws.Cells.Item(50,5).Value = "  -2.55%  "
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "15.91"
$ws.Cells.Item(51,5).Value = "  -5.06%  "
